# Apply weekly update to the Fruta/Tuna subset sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: date moves to 44606, quality/volume/prices shift to "Primera" tier ---
$ws.Range("D6").Value = 44606
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 240
$ws.Range("N6").Value = 11500
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11750
$ws.Range("S6").Value = 653

# --- Row 7: date moves to 44606, quality/prices shift to "Segunda" tier ---
$ws.Range("D7").Value = 44606
$ws.Range("L7").Value = "Segunda"
$ws.Range("N7").Value = 9500
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 9750
$ws.Range("S7").Value = 542

# --- Row 8: date stays 44294, quality/prices shift to "Especial" tier ---
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 14500
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14750
$ws.Range("S8").Value = 819

# --- Row 9: date becomes 44294 (was 44595); values become the old "Primera"/44294 row ---
$ws.Range("D9").Value = 44294
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 240
$ws.Range("N9").Value = 12500
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 12750
$ws.Range("S9").Value = 708

# --- Row 10 (new): the old "Segunda"/44294 record ---
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Terminal La Palmera de La Serena"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44294
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107011
$ws.Range("J10").Value = "Tuna"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 240
$ws.Range("N10").Value = 10500
$ws.Range("O10").Value = 11000
$ws.Range("P10").Value = 10750
$ws.Range("Q10").Value = "$/caja 18 kilos"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 597
$ws.Range("T10").Value = 18

# --- Row 11 (new): the old "Primera"/44595 record, unchanged ---
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Terminal La Palmera de La Serena"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44595
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = "Otros"
$ws.Range("I11").Value = 100107011
$ws.Range("J11").Value = "Tuna"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 15500
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 15750
$ws.Range("Q11").Value = "$/caja 18 kilos"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 875
$ws.Range("T11").Value = 18
